$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab) from "sheet" to "g19.1a"
$ws.Name = "g19.1a"

# Set column B width to best-fit the "Taxa de homicídio doloso" content
# (mirrors double-clicking the column B border to auto-size it)
$ws.Columns.Item(2).AutoFit() | Out-Null

# Select the full column B (B1:B1048576) with B1 as the active cell
$ws.Range("B1:B1048576").Select() | Out-Null
